$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Test One"
$ws.Range("B3").Value = "Just testing things out."
$ws.Range("C3").Value = "test@one.com"
$ws.Range("D3").Value = "Not Urgent"
